$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert 30 new rows at row 780 (pushes the existing "Bleecker Street..."
# and "Letterboxd..." blocks down by 30 rows, to 810-832), mirroring the
# existing "service_top" blocks (apple_top, disney_top, ... vudu_top)
# that already occupy rows 386-779.
$ws.Rows.Item(780).Resize(30).Insert()

$firstRow = 780
$lastRow = 809

# Column F/H/K constant text values (shared by every row in this table).
$fontName = " ComfortAa-Medium"
$white = " #FFFFFF"
$starzColor = " #464646"

for ($i = 0; $i -lt 30; $i++) {
    $r = $firstRow + $i
    $n = $i + 1

    $ws.Cells.Item($r, 1).Value = "TOP_" + $n        # A: key_name
    $ws.Cells.Item($r, 2).Value = "Starz.png"        # B: logo
    $ws.Cells.Item($r, 3).Value = -500                # C: logo_offset
    $ws.Cells.Item($r, 4).Value = 1500                # D: logo_resize
    $ws.Cells.Item($r, 5).Value = 850                 # E: text_offset
    $ws.Cells.Item($r, 6).Value = $fontName           # F: font
    $ws.Cells.Item($r, 8).Value = $white              # H: font_color
    $ws.Cells.Item($r, 9).Value = 0                   # I: border
    $ws.Cells.Item($r, 10).Value = 15                 # J: border_width
    $ws.Cells.Item($r, 11).Value = $white             # K: border_color
    $ws.Cells.Item($r, 14).Value = $starzColor        # N: base_color
    $ws.Cells.Item($r, 15).Value = 1                  # O: gradient
    $ws.Cells.Item($r, 16).Value = 1                  # P: clean
    $ws.Cells.Item($r, 17).Value = 0                  # Q: avg_color
    $ws.Cells.Item($r, 18).Value = " 1',"             # R: white_wash
}

# M: out_name formula -> "starz_" & LOWER(A<row>)
$ws.Range("M780:M809").Formula = '="starz_" & LOWER(A780)'

# S: the big concatenation formula, same pattern used by every other
# block in the sheet (fills relative to each row like Excel's AutoFill).
$sFormula = '="'+"'"+'"&MID(A780,FIND(MID(TRIM(A780),1,1),A780),LEN(A780))&"| "&MID(B780,FIND(MID(TRIM(B780),1,1),B780),LEN(B780))&"| +"&C780&"| "&D780&"| +"&E780&"| "&MID(F780,FIND(MID(TRIM(F780),1,1),F780),LEN(F780))&"| "&G780&"| "&MID(H780,FIND(MID(TRIM(H780),1,1),H780),LEN(H780))&"| "&I780&"| "&J780&"| "&MID(K780,FIND(MID(TRIM(K780),1,1),K780),LEN(K780))&"| "&L780&"| "&MID(M780,FIND(MID(TRIM(M780),1,1),M780),LEN(M780))&"| "&MID(N780,FIND(MID(TRIM(N780),1,1),N780),LEN(N780))&"| "&MID(O780,FIND(MID(TRIM(O780),1,1),O780),LEN(O780))&"| "&MID(P780,FIND(MID(TRIM(P780),1,1),P780),LEN(P780))&"| "&MID(Q780,FIND(MID(TRIM(Q780),1,1),Q780),LEN(Q780))&"| "&MID(R780,FIND(MID(TRIM(R780),1,1),R780),LEN(R780))'
$ws.Range("S780:S809").Formula = $sFormula

$wb.Save()
